$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resolved incorrect budget values on row 4 (B4:E4 were 1, should be 1000)
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 1000
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1000

# Update the active selection to match the saved cursor position
$ws.Range("E14").Select()
